$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    3  = @{ O = "1334"; P = "2330.9k" }
    4  = @{ O = "2";    P = "7324" }
    5  = @{ O = "1980"; P = "1609.6k" }
    6  = @{ O = "1278"; P = "2624.3k" }
    7  = @{ O = "2";    P = "8941" }
    8  = @{ O = "6";    P = "6499" }
    9  = @{ O = "2";    P = "14574" }
    10 = @{ O = "1508"; P = "3126.8k" }
    11 = @{ O = "1924"; P = "1204.6k" }
    12 = @{ O = "3";    P = "7280" }
    13 = @{ O = "2";    P = "4351" }
    14 = @{ O = "3";    P = "2811" }
    15 = @{ O = "4";    P = "10922" }
    16 = @{ O = "1480"; P = "1420.8k" }
    17 = @{ O = "1676"; P = "4202.3k" }
    18 = @{ O = "1631"; P = "2015.4k" }
    19 = @{ O = "1738"; P = "3516.5k" }
    20 = @{ O = "3";    P = "1750" }
    21 = @{ O = "1869"; P = "7415.9k" }
    22 = @{ O = "2";    P = "1636" }
    23 = @{ O = "1636"; P = "2706.6k" }
}

$ws.Range("O3:P23").NumberFormat = "@"

foreach ($row in $values.Keys) {
    $ws.Range("O$row").Value = $values[$row].O
    $ws.Range("P$row").Value = $values[$row].P
}
